$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for serial number 186232438 (row 3) should instead show 186232440
# (the value currently on row 5). Copy it over so the cell keeps its text
# (shared-string) type instead of being re-interpreted as a number.
$ws.Range("B5").Copy($ws.Range("B3"))

# Remove the now-redundant rows (old row 4 "186232437" and old row 5,
# whose value we already moved up into row 3).
$ws.Rows("4:5").Delete()
